# Generate Report for Handoff
#
# The f31b3b96-... file has moved from "Handed back: in sync with en-US" to
# "Ready for handoff" for both locales, with freshly bumped "Latest Handoff"
# timestamps. Update the Overview roll-up sheet plus each per-locale detail
# sheet to match.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for f31b3b96-...md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-32-12 06:32:58"

# --- zh-cn detail sheet: row for f31b3b96-...md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-12 06:32:55"

# --- de-de detail sheet: row for f31b3b96-...md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-12 06:32:58"
